$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round the numeric data in B2:E13 to the nearest integer, as the
# Ontpl_/Pot_ reference files now only need to be written as integer data.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = [double]$cell.Value()
        $rounded = [Math]::Floor($val + 0.5)
        $cell.Value = $rounded
    }
}

$wb.Save()
